# Weekly update: a new week's price record for "Puerro" (Vega Central
# Mapocho de Santiago) is inserted at the top of the data block (row 50),
# pushing the existing rows 50-76 down to 51-77 (the sheet's dimension
# grows from A1:R76 to A1:R77).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 50; this shifts rows 50:76 down
# to 51:77 and carries the existing formatting (incl. the date style on
# column D) along with it.
$ws.Rows.Item(50).Insert()

# Populate the newly inserted row 50 with this week's record.
$ws.Cells.Item(50, 1).Value = 9
$ws.Cells.Item(50, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(50, 3).Value = "Metropolitana"
$ws.Cells.Item(50, 4).Value = 44510
$ws.Cells.Item(50, 5).Value = 13
$ws.Cells.Item(50, 6).Value = 100112005
$ws.Cells.Item(50, 7).Value = "Puerro"
$ws.Cells.Item(50, 8).Value = "Sin especificar"
$ws.Cells.Item(50, 9).Value = "Primera"
$ws.Cells.Item(50, 10).Value = 160
$ws.Cells.Item(50, 11).Value = 6000
$ws.Cells.Item(50, 12).Value = 7000
$ws.Cells.Item(50, 13).Value = 6500
$ws.Cells.Item(50, 14).Value = '$/paquete 20 unidades'
$ws.Cells.Item(50, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(50, 16).Value = 325
$ws.Cells.Item(50, 17).Value = 20
$ws.Cells.Item(50, 18).Value = "Hortaliza"
